# Re-order / rewrite the three Cypher query cells on the "startup" sheet:
#  - B2 (CasesTab row) now holds the Cases query, with the trailing
#    "Cohort" column removed from the RETURN clause.
#  - B3 (SamplesTab row) now holds the Sample query (one blank line removed).
#  - B4 (FilesTab row) keeps the File query (unchanged text).
# Backtick characters (used for Cypher identifier quoting, e.g. `Case ID`)
# are represented with the placeholder "\u00a7" inside the here-strings below,
# because backtick is PowerShell's escape character. They are restored to
# literal backticks afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2Raw = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS §Case ID§ ,
        coalesce(s.clinical_study_designation, '') AS §Study Code§ ,
        coalesce(s.clinical_study_type, '') AS  §Study Type§,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS §Stage of Disease§ ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS §Neutered Status§,
        coalesce(demo.weight, '') AS §Weight (kg)§,
        coalesce(diag.best_response, '') AS §Response to Treatment§
"@
$b3Raw = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS §Sample ID§, 
        coalesce(c.case_id, '') AS §Case ID§, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS §Sample Site§,
        coalesce(samp.summarized_sample_type, '') AS §Sample Type§,
        coalesce(samp.specific_sample_pathology, '') AS §Pathology/Morphology§,
        coalesce(samp.tumor_grade, '') AS §Tumor Grade§,
        coalesce(samp.sample_chronology, '') AS §Sample Chronology§,
        coalesce(samp.percentage_tumor, '') AS §Percentage Tumor§,
        coalesce(samp.necropsy_sample, '') AS §Necropsy Sample§,
        coalesce(samp.sample_preservation, '') AS §Sample Preservation§
"@
$b4Raw = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS §File Name§, 
        coalesce(f.file_type, '') AS §File Type§, 
        coalesce(labels(parent)[0], '') AS §Association§,
        coalesce(f.file_description, '') AS §Description§,
        coalesce(f.file_format, '') AS §File Format§,
        coalesce(f.file_size, '') AS §Size§,
        coalesce(c.case_id, '') AS §Case ID§, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS §Study Code§
"@

$bt = [char]96
$ph = [char]0x00A7
$b2Text = $b2Raw.Replace($ph, $bt)
$b3Text = $b3Raw.Replace($ph, $bt)
$b4Text = $b4Raw.Replace($ph, $bt)

$ws.Range("B2").Value = $b2Text
$ws.Range("B3").Value = $b3Text
$ws.Range("B4").Value = $b4Text

# Row 2 shrinks (one fewer wrapped line) now that the Cohort column was
# dropped from the Cases query text.
$ws.Rows.Item(2).RowHeight = 300

# Move the view/selection from D4 to B2.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("B2").Select()
